$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 150 (previously missing facilidad permanente columns)
$ws.Range("C150").Value = 1
$ws.Range("D150").Value = 0.5

# New daily observations for 06-08-2021 .. 07-09-2021.
# Dates must land as plain text (shared strings), matching the existing
# column A formatting, not auto-converted Excel date serials. We build
# each date value in an off-sheet scratch cell (prefixed with a space so
# Excel's smart-entry parser treats it as text), strip the leading space
# with a formula, then paste the computed text *by value* into column A
# so no numeric/date conversion and no new number-format styles happen.
$dates = @(
    "06-08-2021","09-08-2021","10-08-2021","11-08-2021","12-08-2021","13-08-2021",
    "16-08-2021","17-08-2021","18-08-2021","19-08-2021","20-08-2021","23-08-2021",
    "24-08-2021","25-08-2021","26-08-2021","27-08-2021","30-08-2021","31-08-2021",
    "01-09-2021","02-09-2021","03-09-2021","06-09-2021","07-09-2021"
)
$tpm  = @(0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,0.75,1.5,1.5,1.5,1.5,1.5)
$liq  = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1.75,1.75,1.75,1.75,1.75)
$dep  = @(0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,0.5,1.25,1.25,1.25,1.25,1.25)

$startRow = 151
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    $ws.Range("Z1").Value = " " + $dates[$i]
    $ws.Range("Z2").Formula = "=MID(Z1,2,10)"
    $ws.Range("Z2").Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $tpm[$i]
    $ws.Cells.Item($row, 3).Value = $liq[$i]
    $ws.Cells.Item($row, 4).Value = $dep[$i]
}

$ws.Range("Z1:Z2").Clear()
